$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the authoritative diff: row -> {column -> new value}
$updates = @(
    @{ Row=2; D='62.230.93'; E='  +2.65%  ' }
    @{ Row=3; D='2.425.09'; E='  +3.47%  ' }
    @{ Row=4; E='  +0.21%  ' }
    @{ Row=5; D='555.83'; E='  +2.26%  ' }
    @{ Row=6; D='143.57'; E='  +5.65%  ' }
    @{ Row=7; E='  +0.14%  ' }
    @{ Row=8; D='0.535'; E='  +2.40%  ' }
    @{ Row=9; D='2.426.23'; E='  +3.56%  ' }
    @{ Row=10; E='  +4.96%  ' }
    @{ Row=11; E='  +0.73%  ' }
    @{ Row=12; E='  +2.45%  ' }
    @{ Row=13; E='  +4.32%  ' }
    @{ Row=14; D='26.38'; E='  +7.75%  ' }
    @{ Row=15; E='  +9.83%  ' }
    @{ Row=16; D='2.862.93'; E='  +3.47%  ' }
    @{ Row=17; D='62.104.05'; E='  +2.60%  ' }
    @{ Row=18; D='2.422.07'; E='  +3.35%  ' }
    @{ Row=19; D='11.13'; E='  +5.42%  ' }
    @{ Row=20; D='325.36'; E='  +2.45%  ' }
    @{ Row=21; D='4.19' }
    @{ Row=22; D='6.77'; E='  +3.63%  ' }
    @{ Row=23; E='  +0.27%  ' }
    @{ Row=24; E='  +3.77%  ' }
    @{ Row=25; D='65.00'; E='  +3.15%  ' }
    @{ Row=26; D='9.22'; E='  +11.04%  ' }
    @{ Row=27; D='572.25'; E='  +15.26%  ' }
    @{ Row=28; D='1.00'; E='  +0.11%  ' }
    @{ Row=29; D='2.524.41'; E='  +2.94%  ' }
    @{ Row=30; D='8.40' }
    @{ Row=31; D='0.0₃0939'; E='  +9.88%  ' }
    @{ Row=32; E='  +6.26%  ' }
    @{ Row=33; D='0.149'; E='  +2.97%  ' }
    @{ Row=34; D='1.87'; E='  +5.26%  ' }
    @{ Row=35; E='  +4.61%  ' }
    @{ Row=36; D='5.77'; E='  +11.12%  ' }
    @{ Row=37; B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='4.83'; E='  +6.09%  ' }
    @{ Row=38; E='  +0.03%  ' }
    @{ Row=39; B='Stacks'; C='https://coinranking.com/coin/mMPrMcB7+stacks-stx'; D='1.96'; E='  +9.37%  ' }
    @{ Row=40; E='  +2.83%  ' }
    @{ Row=41; D='18.87'; E='  +2.03%  ' }
    @{ Row=42; D='146.91'; E='  +4.10%  ' }
    @{ Row=43; E='  +0.14%  ' }
    @{ Row=44; D='41.50'; E='  +2.39%  ' }
    @{ Row=45; B='dogwifhat'; C='https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'; D='2.31'; E='  +12.61%  ' }
    @{ Row=46; B='Aave'; C='https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D='151.98'; E='  +7.50%  ' }
    @{ Row=47; D='3.65'; E='  +3.14%  ' }
    @{ Row=48; E='  +7.29%  ' }
    @{ Row=49; D='20.51'; E='  +8.44%  ' }
    @{ Row=50; E='  +4.75%  ' }
    @{ Row=51; D='0.0913'; E='  +1.82%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row

    foreach ($col in @('B', 'C', 'D', 'E')) {
        if (-not $u.ContainsKey($col)) { continue }

        $newValue = $u[$col]
        $cell = $ws.Range("$col$row")

        # Every column in this sheet holds plain text (coin name / link /
        # price / volume), but Excel's COM layer auto-coerces a Value
        # assignment that parses as a plain number (e.g. "555.83") into a
        # numeric cell. Values that are not valid plain numbers (e.g.
        # "62.230.93", which has two dots, or strings padded with
        # spaces/percent signs) are left alone as text automatically. For
        # the ones that would be misread as numbers, force text storage via
        # a transient "@" (text) number format, then restore the cell to
        # the default "Normal" style so no spurious formatting change is
        # left behind.
        $isPlainNumber = $newValue -match '^[0-9]+(\.[0-9]+)?$'

        if ($col -eq 'D' -and $isPlainNumber) {
            $cell.NumberFormat = "@"
            $cell.Value = $newValue
            $cell.Style = "Normal"
        } else {
            $cell.Value = $newValue
        }
    }
}
